# Refactored FlightBookingTest and modified data sheet:
# add "From"/"To" (Bangalore -> Delhi) columns next to the existing
# Locality / Travellers columns on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("C1").Value = "From"
$ws.Range("D1").Value = "To"

# New data row values
$ws.Range("C2").Value = "Bangalore"
$ws.Range("D2").Value = "Delhi"

# Let column C size itself to its new contents (matches the bestFit column
# behaviour already present on column A).
$ws.Columns("C").AutoFit()

# Match the saved selection/active cell state.
$ws.Range("D8").Select()
